$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Val
    )
    $rng = $ws.Range($CellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "61.674.12"
Set-TextValue "E2" "  -3.63%  "
Set-TextValue "D3" "2.479.61"
Set-TextValue "E3" "  -6.17%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "556.71"
Set-TextValue "E5" "  -4.18%  "
Set-TextValue "D6" "147.46"
Set-TextValue "E6" "  -5.82%  "
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "D8" "0.601"
Set-TextValue "E8" "  -4.13%  "
Set-TextValue "D9" "2.477.10"
Set-TextValue "E9" "  -6.20%  "
Set-TextValue "D10" "0.108"
Set-TextValue "E10" "  -8.65%  "
Set-TextValue "D11" "5.49"
Set-TextValue "E11" "  -5.63%  "
Set-TextValue "E12" "  -1.52%  "
Set-TextValue "D13" "0.358"
Set-TextValue "E13" "  -6.87%  "
Set-TextValue "D14" "26.55"
Set-TextValue "E14" "  -7.44%  "
Set-TextValue "D15" "2.926.40"
Set-TextValue "E15" "  -6.20%  "
Set-TextValue "D16" "0.0000168"
Set-TextValue "E16" "  -8.95%  "
Set-TextValue "D17" "61.576.76"
Set-TextValue "E17" "  -3.65%  "
Set-TextValue "D18" "2.483.90"
Set-TextValue "E18" "  -5.94%  "
Set-TextValue "D19" "11.20"
Set-TextValue "E19" "  -8.13%  "
Set-TextValue "D20" "7.20"
Set-TextValue "E20" "  -7.99%  "
Set-TextValue "D21" "4.24"
Set-TextValue "E21" "  -6.54%  "
Set-TextValue "D22" "322.16"
Set-TextValue "E22" "  -6.83%  "
Set-TextValue "E23" "  +0.01%  "
Set-TextValue "D24" "1.90"
Set-TextValue "E24" "  +2.64%  "
Set-TextValue "D25" "64.43"
Set-TextValue "E25" "  -5.59%  "
Set-TextValue "D26" "0.0000100"
Set-TextValue "E26" "  -10.34%  "
Set-TextValue "D27" "564.37"
Set-TextValue "E27" "  -4.05%  "
Set-TextValue "D28" "2.608.92"
Set-TextValue "E28" "  -6.03%  "
Set-TextValue "D29" "1.51"
Set-TextValue "E29" "  -8.18%  "
Set-TextValue "D30" "1.00"
Set-TextValue "E30" "  -0.10%  "
Set-TextValue "D31" "8.38"
Set-TextValue "E31" "  -10.87%  "
Set-TextValue "D32" "7.77"
Set-TextValue "E32" "  -5.45%  "
Set-TextValue "D33" "0.150"
Set-TextValue "E33" "  -6.92%  "
Set-TextValue "E34" "  -6.08%  "
Set-TextValue "E35" "  -7.62%  "
Set-TextValue "D36" "5.97"
Set-TextValue "E36" "  -10.16%  "
Set-TextValue "D37" "4.95"
Set-TextValue "E37" "  -10.00%  "
Set-TextValue "D38" "0.999"
Set-TextValue "E38" "  -0.05%  "
Set-TextValue "D39" "0.385"
Set-TextValue "E39" "  -4.57%  "
Set-TextValue "D40" "18.62"
Set-TextValue "E40" "  -5.90%  "
Set-TextValue "D41" "145.44"
Set-TextValue "E41" "  -4.17%  "
Set-TextValue "E42" "  -7.36%  "
Set-TextValue "E43" "  +0.07%  "
Set-TextValue "E44" "  -4.38%  "
Set-TextValue "D45" "40.63"
Set-TextValue "E45" "  -3.18%  "
Set-TextValue "D46" "148.55"
Set-TextValue "E46" "  -9.25%  "
Set-TextValue "B47" "InjectiveProtocol"
Set-TextValue "C47" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D47" "22.20"
Set-TextValue "E47" "  -9.17%  "
Set-TextValue "B48" "Filecoin"
Set-TextValue "C48" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D48" "3.65"
Set-TextValue "E48" "  -6.70%  "
Set-TextValue "D49" "0.0544"
Set-TextValue "E49" "  -7.84%  "
Set-TextValue "D50" "0.599"
Set-TextValue "E50" "  -5.70%  "
Set-TextValue "D51" "0.0945"
Set-TextValue "E51" "  -5.79%  "
